$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 5142.643
$ws.Cells.Item(62, 9).Value = 3857.8572
$ws.Cells.Item(62, 10).Value = 6427.4287
$ws.Cells.Item(62, 11).Value = 3857.8572
$ws.Cells.Item(62, 12).Value = 6427.4287
$ws.Cells.Item(62, 13).Value = -3233.8572
$ws.Cells.Item(62, 14).Value = -7675.4287
$ws.Cells.Item(65, 8).Value = 5142.643
$ws.Cells.Item(65, 9).Value = 3857.8572
$ws.Cells.Item(65, 10).Value = 6427.4287
$ws.Cells.Item(65, 11).Value = 19289.286
$ws.Cells.Item(65, 12).Value = 32137.1435
$ws.Cells.Item(65, 13).Value = -16169.286
$ws.Cells.Item(65, 14).Value = -38377.14350000001
$ws.Cells.Item(107, 8).Value = 456.9524
$ws.Cells.Item(107, 9).Value = 459.125
$ws.Cells.Item(107, 10).Value = 450
$ws.Cells.Item(107, 11).Value = 459.125
$ws.Cells.Item(107, 12).Value = 450
$ws.Cells.Item(107, 13).Value = 1460.875
$ws.Cells.Item(107, 14).Value = -4290
$ws.Cells.Item(137, 8).Value = 1839.8918
$ws.Cells.Item(137, 9).Value = 1248.4584
$ws.Cells.Item(137, 10).Value = 2931.7693
$ws.Cells.Item(137, 11).Value = 3745.3752
$ws.Cells.Item(137, 12).Value = 8795.3079
$ws.Cells.Item(137, 13).Value = -1195.3752
$ws.Cells.Item(137, 14).Value = -13895.3079
$ws.Cells.Item(138, 8).Value = 2412466.8
$ws.Cells.Item(138, 9).Value = 1218.4166
$ws.Cells.Item(138, 10).Value = 3393313.5
$ws.Cells.Item(138, 11).Value = 3655.2498
$ws.Cells.Item(138, 12).Value = 10179940.5
$ws.Cells.Item(138, 13).Value = 1484.7502
$ws.Cells.Item(138, 14).Value = -10190220.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4829.6616
$ws.Cells.Item(32, 9).Value = 3244.8982
$ws.Cells.Item(32, 11).Value = 3244.8982
$ws.Cells.Item(32, 13).Value = -2957.8982
$ws.Cells.Item(45, 8).Value = 1603.0333
$ws.Cells.Item(45, 9).Value = 1308.6666
$ws.Cells.Item(45, 10).Value = 1729.1904
$ws.Cells.Item(45, 11).Value = 1308.6666
$ws.Cells.Item(45, 12).Value = 1729.1904
$ws.Cells.Item(45, 13).Value = -931.6666
$ws.Cells.Item(45, 14).Value = -2483.1904
$ws.Cells.Item(61, 13).Value = -982.1500000000001
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 1194.15
$ws.Cells.Item(61, 9).Value = 1194.15
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 1194.15
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(74, 8).Value = 56021.3
$ws.Cells.Item(74, 9).Value = 59794.824
$ws.Cells.Item(74, 11).Value = 59794.824
$ws.Cells.Item(74, 13).Value = -58920.824
$ws.Cells.Item(77, 8).Value = 56021.3
$ws.Cells.Item(77, 9).Value = 59794.824
$ws.Cells.Item(77, 11).Value = 298974.12
$ws.Cells.Item(77, 13).Value = -294606.12
$ws.Cells.Item(122, 8).Value = 2064.8
$ws.Cells.Item(122, 9).Value = 1524.6666
$ws.Cells.Item(122, 10).Value = 2875
$ws.Cells.Item(122, 11).Value = 4573.9998
$ws.Cells.Item(122, 12).Value = 8625
$ws.Cells.Item(122, 13).Value = -2123.9998
$ws.Cells.Item(122, 14).Value = -13525
$ws.Cells.Item(132, 8).Value = 2164.5593
$ws.Cells.Item(132, 9).Value = 1985.3959
$ws.Cells.Item(132, 10).Value = 2946.3635
$ws.Cells.Item(132, 11).Value = 5956.1877
$ws.Cells.Item(132, 12).Value = 8839.0905
$ws.Cells.Item(132, 13).Value = -3426.1877
$ws.Cells.Item(132, 14).Value = -13899.0905
$ws.Cells.Item(134, 8).Value = 35000
$ws.Cells.Item(134, 10).Value = 35000
$ws.Cells.Item(134, 12).Value = 35000
$ws.Cells.Item(134, 14).Value = -45140
$ws.Cells.Item(136, 13).Value = -1032.45
$ws.Cells.Item(136, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 1194.15
$ws.Cells.Item(136, 9).Value = 1194.15
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 3582.45
$ws.Cells.Item(136, 12).Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1195.2727
$ws.Cells.Item(94, 9).Value = 1026.8572
$ws.Cells.Item(94, 10).Value = 1490
$ws.Cells.Item(94, 11).Value = 1026.8572
$ws.Cells.Item(94, 12).Value = 1490
$ws.Cells.Item(94, 13).Value = -575.8571999999999
$ws.Cells.Item(94, 14).Value = -2392
$ws.Cells.Item(134, 8).Value = 4671.3477
$ws.Cells.Item(134, 9).Value = 3502.9285
$ws.Cells.Item(134, 11).Value = 10508.7855
$ws.Cells.Item(134, 13).Value = -7973.7855

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 21278192
$ws.Cells.Item(31, 9).Value = 29412702
$ws.Cells.Item(31, 10).Value = 3318.8462
$ws.Cells.Item(31, 11).Value = 29412702
$ws.Cells.Item(31, 12).Value = 3318.8462
$ws.Cells.Item(31, 13).Value = -29412407
$ws.Cells.Item(31, 14).Value = -3908.8462
$ws.Cells.Item(34, 8).Value = 21278192
$ws.Cells.Item(34, 9).Value = 29412702
$ws.Cells.Item(34, 10).Value = 3318.8462
$ws.Cells.Item(34, 11).Value = 29412702
$ws.Cells.Item(34, 12).Value = 3318.8462
$ws.Cells.Item(34, 13).Value = -29412500
$ws.Cells.Item(34, 14).Value = -3722.8462
$ws.Cells.Item(58, 8).Value = 984.8182
$ws.Cells.Item(58, 9).Value = 881.44446
$ws.Cells.Item(58, 10).Value = 1450
$ws.Cells.Item(58, 11).Value = 881.44446
$ws.Cells.Item(58, 12).Value = 1450
$ws.Cells.Item(58, 13).Value = -678.44446
$ws.Cells.Item(58, 14).Value = -1856
$ws.Cells.Item(132, 8).Value = 2232.6667
$ws.Cells.Item(132, 9).Value = 1753.3214
$ws.Cells.Item(132, 10).Value = 3910.375
$ws.Cells.Item(132, 11).Value = 5259.9642
$ws.Cells.Item(132, 12).Value = 11731.125
$ws.Cells.Item(132, 13).Value = -2729.9642
$ws.Cells.Item(132, 14).Value = -16791.125
$ws.Cells.Item(134, 8).Value = 39288030
$ws.Cells.Item(134, 9).Value = 4764561
$ws.Cells.Item(134, 10).Value = 142858450
$ws.Cells.Item(134, 11).Value = 14293683
$ws.Cells.Item(134, 12).Value = 428575350
$ws.Cells.Item(134, 13).Value = -14291148
$ws.Cells.Item(134, 14).Value = -428580420
$ws.Cells.Item(136, 8).Value = 984.8182
$ws.Cells.Item(136, 9).Value = 881.44446
$ws.Cells.Item(136, 10).Value = 1450
$ws.Cells.Item(136, 11).Value = 2644.33338
$ws.Cells.Item(136, 12).Value = 4350
$ws.Cells.Item(136, 13).Value = -94.33338000000003
$ws.Cells.Item(136, 14).Value = -9450

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 14103362
$ws.Cells.Item(5, 9).Value = 18182310
$ws.Cells.Item(5, 10).Value = 11112134
$ws.Cells.Item(5, 11).Value = 54546930
$ws.Cells.Item(5, 12).Value = 33336402
$ws.Cells.Item(5, 13).Value = -54546818
$ws.Cells.Item(5, 14).Value = -33336626
$ws.Cells.Item(122, 8).Value = 934.8461
$ws.Cells.Item(122, 9).Value = 536.7059
$ws.Cells.Item(122, 10).Value = 1686.8889
$ws.Cells.Item(122, 11).Value = 4830.3531
$ws.Cells.Item(122, 12).Value = 15182.0001
$ws.Cells.Item(122, 13).Value = -2380.3531
$ws.Cells.Item(122, 14).Value = -20082.0001
$ws.Cells.Item(131, 8).Value = 889.97
$ws.Cells.Item(131, 9).Value = 612.7143
$ws.Cells.Item(131, 10).Value = 910.8387
$ws.Cells.Item(131, 11).Value = 1838.1429
$ws.Cells.Item(131, 12).Value = 2732.5161
$ws.Cells.Item(131, 13).Value = 3201.8571
$ws.Cells.Item(131, 14).Value = -12812.5161
$ws.Cells.Item(134, 8).Value = 7533
$ws.Cells.Item(134, 9).Value = 6821.5
$ws.Cells.Item(134, 10).Value = 8102.2
$ws.Cells.Item(134, 11).Value = 20464.5
$ws.Cells.Item(134, 12).Value = 24306.6
$ws.Cells.Item(134, 13).Value = -15394.5
$ws.Cells.Item(134, 14).Value = -34446.6
$ws.Cells.Item(135, 8).Value = 14103362
$ws.Cells.Item(135, 9).Value = 18182310
$ws.Cells.Item(135, 10).Value = 11112134
$ws.Cells.Item(135, 11).Value = 163640790
$ws.Cells.Item(135, 12).Value = 100009206
$ws.Cells.Item(135, 13).Value = -163638255
$ws.Cells.Item(135, 14).Value = -100014276

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4092.04
$ws.Cells.Item(132, 9).Value = 4064.3684
$ws.Cells.Item(132, 10).Value = 4179.6665
$ws.Cells.Item(132, 11).Value = 12193.1052
$ws.Cells.Item(132, 12).Value = 12538.9995
$ws.Cells.Item(132, 13).Value = -9663.1052
$ws.Cells.Item(132, 14).Value = -17598.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(134, 8).Value = 32312
$ws.Cells.Item(134, 10).Value = 32312
$ws.Cells.Item(134, 12).Value = 32312
$ws.Cells.Item(134, 14).Value = -42452
$ws.Cells.Item(136, 8).Value = 12821980
$ws.Cells.Item(136, 9).Value = 13334779
$ws.Cells.Item(136, 10).Value = 2000
$ws.Cells.Item(136, 11).Value = 40004337
$ws.Cells.Item(136, 12).Value = 6000
$ws.Cells.Item(136, 13).Value = -40001787
$ws.Cells.Item(136, 14).Value = -11100
$ws.Cells.Item(137, 8).Value = 46060
$ws.Cells.Item(137, 9).Value = 28593.334
$ws.Cells.Item(137, 10).Value = 51300
$ws.Cells.Item(137, 11).Value = 28593.334
$ws.Cells.Item(137, 13).Value = -23493.334
$ws.Cells.Item(137, 14).Value = -61500

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(86, 8).Value = 29800
$ws.Cells.Item(86, 10).Value = 29800
$ws.Cells.Item(86, 12).Value = 29800
$ws.Cells.Item(86, 14).Value = -32046
$ws.Cells.Item(89, 8).Value = 29800
$ws.Cells.Item(89, 10).Value = 29800
$ws.Cells.Item(89, 12).Value = 149000
$ws.Cells.Item(89, 14).Value = -160232
$ws.Cells.Item(122, 8).Value = 10475980
$ws.Cells.Item(122, 9).Value = 13900617
$ws.Cells.Item(122, 10).Value = 202067.5
$ws.Cells.Item(122, 11).Value = 41701851
$ws.Cells.Item(122, 12).Value = 606202.5
$ws.Cells.Item(122, 13).Value = -41699401
$ws.Cells.Item(122, 14).Value = -611102.5
$ws.Cells.Item(132, 8).Value = 4169285.5
$ws.Cells.Item(132, 9).Value = 4547914
$ws.Cells.Item(132, 10).Value = 4374.5
$ws.Cells.Item(132, 11).Value = 13643742
$ws.Cells.Item(132, 12).Value = 13123.5
$ws.Cells.Item(132, 13).Value = -13641212
$ws.Cells.Item(132, 14).Value = -18183.5
$ws.Cells.Item(136, 8).Value = 22729616
$ws.Cells.Item(136, 9).Value = 33335814
$ws.Cells.Item(136, 10).Value = 2047.0714
$ws.Cells.Item(136, 11).Value = 100007442
$ws.Cells.Item(136, 12).Value = 6141.2142
$ws.Cells.Item(136, 13).Value = -100004892
$ws.Cells.Item(136, 14).Value = -11241.2142
